# Edit script: add "PO Forecast" sheet, rename header columns on the
# existing "Weekly Quantity" / "Monthly Trend" sheets.

$wb = $excel.ActiveWorkbook

# --- 1) Rename the "Requested quantity" headers -------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after "Monthly Trend" -----------
# NOTE: Worksheets.Add() inserts the new sheet at index 0 and this shifts
# every other sheet's index up by one, which invalidates any worksheet
# object captured *before* the call (the host resolves worksheet objects
# by index, not by a stable handle). So re-fetch everything we need by
# name *after* Add() has run.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "PO Forecast"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$ws = $wb.Worksheets.Item("PO Forecast")
$ws.Move($null, $wsMonthly)

# Re-fetch once more: Move() reorders sheets too, so refresh the handle.
$ws = $wb.Worksheets.Item("PO Forecast")

# --- 3) Header row --------------------------------------------------------
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- 4) Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) --------------
$data = @(
    @(45186.99999999999, 117, -107.8156053361354, 352.8878526359118),
    @(45193.99999999999, 120, -127.9553564065731, 360.4836053134239),
    @(45200.99999999999, 122, -107.0880532916763, 353.1876072324823),
    @(45207.99999999999, 125, -116.8216850841331, 356.7446756697827),
    @(45214.99999999999, 128, -96.87336215715042, 371.4575069693952),
    @(45221.99999999999, 130, -110.8060056995178, 365.1630307527028),
    @(45235.99999999999, 136, -102.4106441473369, 362.5867914373903),
    @(45242.99999999999, 139, -113.60474088651, 373.0090056052965),
    @(45256.99999999999, 144, -93.72028888070989, 375.5970089838463),
    @(45270.99999999999, 149, -87.85748087054564, 382.37638891834),
    @(45305.99999999999, 163, -77.38738791736439, 400.2629960775457),
    @(45312.99999999999, 166, -62.43120816204673, 395.476202352789),
    @(45326.99999999999, 171, -61.51713509115232, 414.9496011690355),
    @(45333.99999999999, 174, -54.08958119714394, 418.8718452218005),
    @(45340.99999999999, 177, -57.32305804885316, 420.8421125291805),
    @(45347.99999999999, 179, -57.17297093155749, 417.5161262756209),
    @(45354.99999999999, 182, -49.02970917071449, 411.6974304729098),
    @(45361.99999999999, 185, -51.8154312328558, 418.6991387805896),
    @(45368.99999999999, 188, -42.38573500454984, 417.9707423723924),
    @(45375.99999999999, 190, -44.5088717253438, 439.5366457072464),
    @(45382.99999999999, 193, -34.7310899184376, 419.1336562286577),
    @(45389.99999999999, 196, -38.01556362398571, 438.7405909714103),
    @(45396.99999999999, 199, -36.0085060271548, 427.9817112156763),
    @(45403.99999999999, 201, -49.95994729018464, 439.4064749490021),
    @(45410.99999999999, 204, -40.56858916520811, 432.3687288136867),
    @(45417.99999999999, 207, -24.06466678379608, 432.5019380378122),
    @(45424.99999999999, 209, -19.81298383345885, 431.9542729503685),
    @(45431.99999999999, 212, -34.88237528440877, 458.0514456595895),
    @(45459.99999999999, 223, -34.25325023612307, 458.9536195565209),
    @(45466.99999999999, 226, 0.9968247477412029, 461.6695466160965),
    @(45473.99999999999, 229, -18.2350422890955, 461.1425696017527),
    @(45480.99999999999, 231, -1.643392889149597, 476.3093084494071),
    @(45487.99999999999, 234, 12.37181936150104, 470.4543329245977),
    @(45494.99999999999, 237, -5.922746428116626, 474.3152894898026),
    @(45515.99999999999, 245, 13.19799350058734, 489.9888389734171),
    @(45529.99999999999, 250, 8.76704516098094, 491.1646070263494),
    @(45536.99999999999, 253, 28.98448228462037, 506.2623726520646),
    @(45557.99999999999, 261, 25.17344832978434, 501.6410652243907),
    @(45564.99999999999, 264, 55.0736751118088, 502.0947779374379),
    @(45571.99999999999, 267, 23.91496440748239, 496.0279620946152),
    @(45585.99999999999, 272, 37.28652609848183, 514.3495320984071),
    @(45592.99999999999, 275, 64.18303946428389, 513.4835015579508),
    @(45599.99999999999, 278, 46.5186694792175, 522.3153836468186),
    @(45606.99999999999, 280, 49.06415393643122, 528.4333495805208),
    @(45613.99999999999, 283, 48.37889144243586, 516.1232701564196),
    @(45620.99999999999, 286, 31.85919456531957, 507.500683835346),
    @(45627.99999999999, 288, 59.85217860457219, 530.457692792282),
    @(45634.99999999999, 291, 52.43309621758998, 517.3712790853696),
    @(45641.99999999999, 294, 57.64004566347948, 534.777170274275),
    @(45648.99999999999, 297, 53.19669807961105, 536.4971219345314),
    @(45655.99999999999, 299, 45.64139261495534, 525.4291754058731)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}

Write-Output "PO Forecast sheet populated."
